$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_OpenBrowser")

# Bring formatting (style s="4", thin borders, wrap text) for the new row
# by copying the format from the row above (row 4 is plain style 4 across
# every column, including the normally-blank E:J cells).
$ws.Range("A4:O4").Copy()
$ws.Range("A5:O5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New test-case row: "Verify Broken Link" / VerifyBrokenLink action.
$ws.Range("A5").Value = "US_ID_01"
$ws.Range("B5").Value = "TC_CityMarket_04"
$ws.Range("C5").Value = "Verify Broken Link"
$ws.Range("D5").Value = "VerifyBrokenLink"
$ws.Range("K5").Value = "All links should workd"
$ws.Range("L5").Value = "VerifyBrokenLink: null"
$ws.Range("M5").Value = "Fail"
$ws.Range("N5").Value = "-"
$ws.Range("O5").Value = "OpenBrowser"

# Match the row height used by the other data rows in this table.
$ws.Rows.Item(5).RowHeight = 30

# Reflect the new active cell / scroll position used while editing.
$ws.Activate() | Out-Null
$ws.Range("N5").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
